# Revert the C7 coefficient back to its earlier value and move the
# active cell selection to C8 (matching the pre-merge worksheet state).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C7").Value = 12.08

$ws.Range("C8").Select()
